$d = $word.ActiveDocument

# --- Mark the title block (paragraph 1) and the following 3 blank paragraphs
# (2-4) as English (US) proofing language, matching the rest of the
# document's headings which already carry lang="en-US". For paragraphs
# without runs this stamps the paragraph-mark run properties (w:pPr/w:rPr);
# for the title paragraph it stamps each run's properties.
foreach ($i in 1, 2, 3, 4) {
    $d.Paragraphs($i).Range.LanguageID = "en-US"
}

# --- Content corrections (names fixed to match the cross-referenced use case) ---

# "UC02 Delete User" -> "UC02 Slet Bruger"
$d.Content.Find.Execute("Delete User", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Slet Bruger", 2)

# "Klienten K er logget ind på brugeren" -> "Patienten P er logget ind på brugeren"
$d.Content.Find.Execute("Klienten K er logget ind på brugeren", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Patienten P er logget ind på brugeren", 2)

# "Anmodning om password blev præsenteret for K" -> "...for P"
$d.Content.Find.Execute("Anmodning om password blev præsenteret for K", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Anmodning om password blev præsenteret for P", 2)
